$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking values (e.g. "0.9992", "9.008")
# are not auto-converted to numbers by Excel, matching the inlineStr type in the diff.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.279.86'
$ws.Range("E2").Value = '  -0.53%  '

$ws.Range("D3").Value = '1.843.64'
$ws.Range("E3").Value = '  -0.35%  '

$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '240.69'
$ws.Range("E5").Value = '  +0.25%  '

$ws.Range("D6").Value = '0.6281'

$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").Value = '0.07448'
$ws.Range("E8").Value = '  -2.68%  '

$ws.Range("D9").Value = '0.2893'
$ws.Range("E9").Value = '  -0.77%  '

$ws.Range("D10").Value = '24.26'
$ws.Range("E10").Value = '  -2.19%  '

$ws.Range("D11").Value = '0.07727'
$ws.Range("E11").Value = '  -0.09%  '

$ws.Range("D12").Value = '1.842.95'
$ws.Range("E12").Value = '  -2.40%  '

$ws.Range("D13").Value = '4.987'
$ws.Range("E13").Value = '  -0.89%  '

$ws.Range("D14").Value = '0.6772'
$ws.Range("E14").Value = '  -0.54%  '

$ws.Range("D15").Value = '0.00001007'
$ws.Range("E15").Value = '  -4.75%  '

$ws.Range("D16").Value = '82.04'
$ws.Range("E16").Value = '  -1.79%  '

$ws.Range("D17").Value = '6.134'
$ws.Range("E17").Value = '  -1.05%  '

$ws.Range("D18").Value = '29.285.98'
$ws.Range("E18").Value = '  -0.98%  '

$ws.Range("D19").Value = '227.83'
$ws.Range("E19").Value = '  -0.56%  '

$ws.Range("D20").Value = '12.27'
$ws.Range("E20").Value = '  -0.54%  '

$ws.Range("E21").Value = '  +0.11%  '

$ws.Range("D22").Value = '7.376'
$ws.Range("E22").Value = '  -1.14%  '

$ws.Range("D23").Value = '0.9995'
$ws.Range("E23").Value = '  -0.12%  '

$ws.Range("D24").Value = '158.78'
$ws.Range("E24").Value = '  +0.96%  '

$ws.Range("D25").Value = '0.1373'
$ws.Range("E25").Value = '  -0.89%  '

$ws.Range("D26").Value = '8.385'
$ws.Range("E26").Value = '  -0.47%  '

$ws.Range("D27").Value = '17.56'
$ws.Range("E27").Value = '  -0.93%  '

$ws.Range("D28").Value = '0.06242'
$ws.Range("E28").Value = '  +11.23%  '

$ws.Range("D29").Value = '1.392'
$ws.Range("E29").Value = '  +0.95%  '

$ws.Range("D30").Value = '1.474'
$ws.Range("E30").Value = '  +0.91%  '

$ws.Range("D31").Value = '4.079'
$ws.Range("E31").Value = '  -1.25%  '

$ws.Range("D32").Value = '4.064'
$ws.Range("E32").Value = '  -0.01%  '

$ws.Range("D33").Value = '1.817'

$ws.Range("D34").Value = '1.140'
$ws.Range("E34").Value = '  -2.28%  '

$ws.Range("D35").Value = '0.6956'
$ws.Range("E35").Value = '  -0.65%  '

$ws.Range("D36").Value = '2.586'
$ws.Range("E36").Value = '  -0.21%  '

$ws.Range("D37").Value = '2.841'
$ws.Range("E37").Value = '  +3.63%  '

$ws.Range("D38").Value = '1.247.95'
$ws.Range("E38").Value = '  +1.67%  '

$ws.Range("D39").Value = '0.01814'
$ws.Range("E39").Value = '  +0.54%  '

$ws.Range("D40").Value = '6.513'
$ws.Range("E40").Value = '  +0.73%  '

$ws.Range("D41").Value = '0.9092'
$ws.Range("E41").Value = '  +0.32%  '

$ws.Range("D42").Value = '0.9993'
$ws.Range("E42").Value = '  -0.10%  '

$ws.Range("D43").Value = '2.000.62'
$ws.Range("E43").Value = '  -15.54%  '

$ws.Range("D44").Value = '101.40'
$ws.Range("E44").Value = '  -0.56%  '

$ws.Range("D45").Value = '66.15'
$ws.Range("E45").Value = '  +0.15%  '

$ws.Range("D46").Value = '7.055'
$ws.Range("E46").Value = '  -2.13%  '

$ws.Range("D47").Value = '0.1163'
$ws.Range("E47").Value = '  +0.81%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '9.008'
$ws.Range("E48").Value = '  +0.26%  '

$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D49").Value = '0.3939'
$ws.Range("E49").Value = '  -2.24%  '

$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.00000000114'
$ws.Range("E50").Value = '  -4.60%  '

$ws.Range("D51").Value = '1.658'
$ws.Range("E51").Value = '  -1.28%  '
